$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price/volume figures refreshed, rows 20/21 and 42/43 reordered)
$updates = @{
    'D2' = '25.878.66'
    'E2' = '  +0.54%  '
    'D3' = '1.630.34'
    'E3' = '  -0.03%  '
    'D4' = '0.998'
    'E4' = '  -0.23%  '
    'D5' = '214.77'
    'E5' = '  +0.17%  '
    'E6' = '  +0.31%  '
    'E7' = '  -0.04%  '
    'E8' = '  +0.10%  '
    'D9' = '0.0631'
    'E9' = '  -0.24%  '
    'D10' = '19.63'
    'E10' = '  +0.64%  '
    'E11' = '  -0.61%  '
    'D12' = '1.858.74'
    'E12' = '  +0.19%  '
    'E13' = '  -0.61%  '
    'D14' = '1.623.96'
    'E14' = '  +0.01%  '
    'D15' = '0.544'
    'E15' = '  -2.15%  '
    'D16' = '0.0₃0757'
    'E16' = '  -0.39%  '
    'D17' = '62.79'
    'E17' = '  -0.17%  '
    'D18' = '25.920.31'
    'E18' = '  +0.80%  '
    'D19' = '0.999'
    'E19' = '  -0.20%  '
    'B20' = 'BitcoinCash'
    'C20' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D20' = '192.81'
    'E20' = '  +0.52%  '
    'B21' = 'Uniswap'
    'C21' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D21' = '4.37'
    'E21' = '  -1.47%  '
    'D22' = '9.94'
    'E22' = '  +0.21%  '
    'E23' = '  +0.19%  '
    'D24' = '1.80'
    'E24' = '  -1.20%  '
    'E25' = '  -0.20%  '
    'D26' = '142.31'
    'E26' = '  +0.07%  '
    'E27' = '  +2.10%  '
    'E28' = '  +0.02%  '
    'D29' = '15.45'
    'E29' = '  -0.09%  '
    'E30' = '  +0.36%  '
    'D31' = '0.0499'
    'E31' = '  +1.99%  '
    'E32' = '  -0.60%  '
    'E33' = '  -0.43%  '
    'E34' = '  -0.17%  '
    'D35' = '2.41'
    'E35' = '  +0.65%  '
    'D36' = '0.901'
    'E36' = '  -0.39%  '
    'D37' = '1.134.50'
    'E37' = '  -0.38%  '
    'D38' = '0.549'
    'E38' = '  +1.24%  '
    'D39' = '2.45'
    'E39' = '  -2.42%  '
    'E40' = '  +0.55%  '
    'D41' = '0.998'
    'E41' = '  -0.10%  '
    'B42' = 'TrustWalletToken'
    'C42' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D42' = '0.803'
    'E42' = '  +0.11%  '
    'B43' = 'FraxShare'
    'C43' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D43' = '5.47'
    'E43' = '  -1.24%  '
    'D44' = '99.11'
    'E44' = '  -1.28%  '
    'D45' = '1.767.92'
    'E45' = '  +0.18%  '
    'E46' = '  +1.89%  '
    'D47' = '56.02'
    'E47' = '  +1.38%  '
    'D48' = '0.0527'
    'E48' = '  +3.90%  '
    'D49' = '1.45'
    'E49' = '  +1.37%  '
    'D50' = '0.413'
    'E50' = '  -1.31%  '
    'D51' = '7.59'
    'E51' = '  +2.51%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so numeric-looking strings (e.g. "214.77") are not
    # coerced into floating point numbers, then strip the temporary formatting so the
    # cell style matches the original (unstyled) cells.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
